# The user can now choose not to save search results, and the Yahoo
# company-name lookup no longer throws unhandled on failure. This run's
# search was saved, appending a new result row to the named-entity log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42602.583043981482
$ws.Range("B5").Value = "Named"
$ws.Range("C5").Value = 8988
$ws.Range("D5").Value = 6493
$ws.Range("E5").Value = 408
$ws.Range("F5").Value = 58
$ws.Range("G5").Value = 27
$ws.Range("H5").Value = 67
$ws.Range("I5").Value = 31
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 0
